$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (27) into two
# new rows, preserving styles/fonts without touching styles.xml.
$ws.Rows("27:27").Copy()
$ws.Rows("28:28").Insert()
$ws.Rows("27:27").Copy()
$ws.Rows("29:29").Insert()

# Row 28
$ws.Range("A28").Value = "com.hamxa.shaynachim"
$ws.Range("B28").Value = "bitcoin"
$ws.Range("C28").Value = "nachumella625@gmail.com"
$ws.Range("D28").Value = "milleradir327@gmail.com "
$ws.Range("E28").Value = "27/5/2019 15:59"
$ws.Range("F28").Value = "very valuable information in this great app"
$ws.Range("G28").Value = "no"

# Row 29
$ws.Range("A29").Value = "com.hamxa.shaynachim"
$ws.Range("B29").Value = "bitcoin"
$ws.Range("C29").Value = "dan624655@gmail.com"
$ws.Range("D29").Value = "nachumella625@gmail.com"
$ws.Range("E29").Value = "27/5/2019 15:59"
$ws.Range("F29").Value = "fantastic info! Guaranteed!"
$ws.Range("G29").Value = "no"

$ws.Application.CutCopyMode = $false

$ws.Range("G30").Select()
